# Username change for overage
#
# The "Input" sheet of ManageProducts.xlsx lists products in column B
# (ProductName). Three existing rows get a new generated product name
# (same "prodXXXXXXXX" naming convention used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "prodhUwz"
$ws.Range("B3").Value = "prodgyDO"
$ws.Range("B5").Value = "prodUDld"
